# Auto-generated edit script applying numeric updates described by the commit diff
# (market-data refresh for the Goblin_Profits workbook; 8 class sheets: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2999.5
$ws.Range("J28").Value = 2999.5
$ws.Range("L28").Value = 2999.5
$ws.Range("N28").Value = -3969.5
$ws.Range("H62").Value = 97523.78
$ws.Range("I62").Value = 206142.75
$ws.Range("J62").Value = 10628.6
$ws.Range("K62").Value = 206142.75
$ws.Range("L62").Value = 10628.6
$ws.Range("M62").Value = -205518.75
$ws.Range("N62").Value = -11876.6
$ws.Range("H65").Value = 97523.78
$ws.Range("I65").Value = 206142.75
$ws.Range("J65").Value = 10628.6
$ws.Range("K65").Value = 1030713.75
$ws.Range("L65").Value = 53143
$ws.Range("M65").Value = -1027593.75
$ws.Range("N65").Value = -59383
$ws.Range("H94").Value = 4959.364
$ws.Range("I94").Value = 4959.364
$ws.Range("K94").Value = 4959.364
$ws.Range("M94").Value = -4508.364
$ws.Range("H106").Value = 2455.1304
$ws.Range("I106").Value = 2225.9
$ws.Range("J106").Value = 3983.3333
$ws.Range("K106").Value = 2225.9
$ws.Range("L106").Value = 3983.3333
$ws.Range("M106").Value = -1594.9
$ws.Range("N106").Value = -5245.3333
$ws.Range("H137").Value = 2724.5217
$ws.Range("I137").Value = 2352.9092
$ws.Range("J137").Value = 3065.1667
$ws.Range("K137").Value = 7058.7276
$ws.Range("L137").Value = 9195.500100000001
$ws.Range("M137").Value = -4508.7276
$ws.Range("N137").Value = -14295.5001
$ws.Range("H138").Value = 1312516.2
$ws.Range("J138").Value = 1591732
$ws.Range("L138").Value = 4775196
$ws.Range("N138").Value = -4785476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 4352.615
$ws.Range("I41").Value = 1229.3334
$ws.Range("K41").Value = 1229.3334
$ws.Range("M41").Value = -815.3334
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988
$ws.Range("H110").Value = 1521.7273
$ws.Range("I110").Value = 1474
$ws.Range("K110").Value = 1474
$ws.Range("M110").Value = 571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 15995
$ws.Range("J46").Value = 15995
$ws.Range("L46").Value = 15995
$ws.Range("N46").Value = -16591
$ws.Range("H107").Value = 4661.8945
$ws.Range("I107").Value = 3686.7144
$ws.Range("J107").Value = 7392.4
$ws.Range("K107").Value = 3686.7144
$ws.Range("L107").Value = 7392.4
$ws.Range("M107").Value = -1766.7144
$ws.Range("N107").Value = -11232.4
$ws.Range("H134").Value = 2648.158
$ws.Range("I134").Value = 2648.158
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7944.474
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5409.474
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 192.64706
$ws.Range("I7").Value = 243.18182
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 243.18182
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = -130.18182
$ws.Range("N7").Value = -326
$ws.Range("H22").Value = 1998.8334
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 2198.8
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 2198.8
$ws.Range("M22").Value = -649
$ws.Range("N22").Value = -2898.8
$ws.Range("H68").Value = 55499.832
$ws.Range("I68").Value = 51599.8
$ws.Range("K68").Value = 51599.8
$ws.Range("M68").Value = -50850.8
$ws.Range("H71").Value = 55499.832
$ws.Range("I71").Value = 51599.8
$ws.Range("K71").Value = 154799.4
$ws.Range("M71").Value = -151055.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1134.5
$ws.Range("J18").Value = 850
$ws.Range("L18").Value = 2550
$ws.Range("N18").Value = -2888
$ws.Range("H102").Value = 40244.6
$ws.Range("I102").Value = 1999
$ws.Range("J102").Value = 49806
$ws.Range("K102").Value = 5997
$ws.Range("L102").Value = 149418
$ws.Range("M102").Value = -3563
$ws.Range("N102").Value = -154286
$ws.Range("H107").Value = 1045.3846
$ws.Range("J107").Value = 864.14703
$ws.Range("L107").Value = 2592.44109
$ws.Range("N107").Value = -6432.44109
$ws.Range("H122").Value = 3946.5
$ws.Range("J122").Value = 4339.875
$ws.Range("L122").Value = 39058.875
$ws.Range("N122").Value = -43958.875
$ws.Range("H129").Value = 4978.9375
$ws.Range("J129").Value = 7205.5
$ws.Range("L129").Value = 21616.5
$ws.Range("N129").Value = -31616.5
$ws.Range("H131").Value = 5612.1924
$ws.Range("I131").Value = 988.8889
$ws.Range("J131").Value = 8059.8237
$ws.Range("K131").Value = 2966.6667
$ws.Range("L131").Value = 24179.4711
$ws.Range("M131").Value = 2073.3333
$ws.Range("N131").Value = -34259.4711
$ws.Range("H133").Value = 10211.462
$ws.Range("I133").Value = 2549.8
$ws.Range("K133").Value = 7649.400000000001
$ws.Range("M133").Value = -2589.400000000001
$ws.Range("H134").Value = 4884.1514
$ws.Range("I134").Value = 1181.4166
$ws.Range("K134").Value = 3544.2498
$ws.Range("M134").Value = 1525.7502
$ws.Range("H138").Value = 2741.75
$ws.Range("J138").Value = 2498.3125
$ws.Range("L138").Value = 7494.9375
$ws.Range("N138").Value = -17774.9375
$ws.Range("H140").Value = 1746.6666
$ws.Range("I140").Value = 1694.48
$ws.Range("K140").Value = 5083.440000000001
$ws.Range("M140").Value = 96.55999999999949

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 3000
$ws.Range("J25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -4058
$ws.Range("H105").Value = 70000
$ws.Range("J105").Value = 70000
$ws.Range("L105").Value = 70000
$ws.Range("N105").Value = -76988
$ws.Range("H107").Value = 1758
$ws.Range("I107").Value = 649.3333
$ws.Range("J107").Value = 2589.5
$ws.Range("K107").Value = 649.3333
$ws.Range("L107").Value = 2589.5
$ws.Range("M107").Value = 1270.6667
$ws.Range("N107").Value = -6429.5
$ws.Range("H126").Value = 4305
$ws.Range("I126").Value = 2687.25
$ws.Range("K126").Value = 8061.75
$ws.Range("M126").Value = -5591.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("H46").Value = 4061.0833
$ws.Range("J46").Value = 4040
$ws.Range("L46").Value = 4040
$ws.Range("N46").Value = -4416
$ws.Range("H68").Value = 6216
$ws.Range("I68").Value = 3448.8
$ws.Range("J68").Value = 7599.6
$ws.Range("K68").Value = 3448.8
$ws.Range("L68").Value = 7599.6
$ws.Range("M68").Value = -2699.8
$ws.Range("N68").Value = -9097.6
$ws.Range("H71").Value = 6216
$ws.Range("I71").Value = 3448.8
$ws.Range("J71").Value = 7599.6
$ws.Range("K71").Value = 17244
$ws.Range("L71").Value = 37998
$ws.Range("M71").Value = -13500
$ws.Range("N71").Value = -45486
$ws.Range("H87").Value = 1000000000
$ws.Range("J87").Value = 1000000000
$ws.Range("L87").Value = 1000000000
$ws.Range("N87").Value = -1000002246
$ws.Range("H90").Value = 1000000000
$ws.Range("J90").Value = 1000000000
$ws.Range("L90").Value = 3000000000
$ws.Range("N90").Value = -3000011232
$ws.Range("H93").Value = 4776.1
$ws.Range("I93").Value = 2663.7058
$ws.Range("J93").Value = 7538.4614
$ws.Range("K93").Value = 2663.7058
$ws.Range("L93").Value = 7538.4614
$ws.Range("M93").Value = -1415.7058
$ws.Range("N93").Value = -10034.4614
$ws.Range("H122").Value = 6250.839
$ws.Range("I122").Value = 6395.3477
$ws.Range("K122").Value = 19186.0431
$ws.Range("M122").Value = -16736.0431
$ws.Range("H127").Value = 200000
$ws.Range("J127").Value = 200000
$ws.Range("L127").Value = 200000
$ws.Range("N127").Value = -209920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9916.583000000001
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 12374.875
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 12374.875
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -13622.875
$ws.Range("H65").Value = 9916.583000000001
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 12374.875
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 61874.375
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -68114.375
$ws.Range("H96").Value = 5966.222
$ws.Range("I96").Value = 5739.2
$ws.Range("K96").Value = 5739.2
$ws.Range("M96").Value = -4366.2
$ws.Range("H100").Value = 840.8823
$ws.Range("I100").Value = 522.7692
$ws.Range("J100").Value = 1874.75
$ws.Range("K100").Value = 1045.5384
$ws.Range("L100").Value = 3749.5
$ws.Range("M100").Value = -504.5383999999999
$ws.Range("N100").Value = -4831.5
$ws.Range("H126").Value = 1915.7858
$ws.Range("I126").Value = 1870.8462
$ws.Range("K126").Value = 5612.5386
$ws.Range("M126").Value = -3142.5386

